{"js": "// Update the worksheet date and the twenty-five division problems.\n// Listed in document order; `418\u00f77=` occurs twice with different\n// replacement values, so replacements are matched up positionally\n// (the Nth occurrence of a given search string maps to the Nth entry\n// below that shares that search string).\nconst replacements = [\n  [\"2024-06-03 Monday\", \"2024-06-04 Tuesday\"],\n  [\"274\u00f76=\", \"687\u00f78=\"],\n  [\"187\u00f76=\", \"305\u00f75=\"],\n  [\"418\u00f77=\", \"244\u00f72=\"],\n  [\"700\u00f78=\", \"930\u00f73=\"],\n  [\"780\u00f77=\", \"846\u00f75=\"],\n  [\"343\u00f73=\", \"134\u00f73=\"],\n  [\"336\u00f77=\", \"944\u00f72=\"],\n  [\"418\u00f77=\", \"708\u00f73=\"],\n  [\"655\u00f74=\", \"245\u00f79=\"],\n  [\"290\u00f75=\", \"502\u00f74=\"],\n  [\"682\u00f76=\", \"923\u00f72=\"],\n  [\"544\u00f77=\", \"606\u00f79=\"],\n  [\"590\u00f73=\", \"768\u00f76=\"],\n  [\"679\u00f74=\", \"160\u00f77=\"],\n  [\"262\u00f79=\", \"109\u00f74=\"],\n  [\"821\u00f73=\", \"999\u00f79=\"],\n  [\"119\u00f75=\", \"755\u00f77=\"],\n  [\"881\u00f77=\", \"740\u00f78=\"],\n  [\"652\u00f78=\", \"281\u00f79=\"],\n  [\"957\u00f73=\", \"183\u00f78=\"],\n  [\"678\u00f74=\", \"540\u00f75=\"],\n  [\"395\u00f76=\", \"898\u00f72=\"],\n  [\"515\u00f78=\", \"768\u00f79=\"],\n  [\"535\u00f75=\", \"857\u00f76=\"],\n  [\"625\u00f73=\", \"308\u00f78=\"],\n];\n\nconst body = context.document.body;\n\n// Search once per distinct \"find\" string, then walk the hits in\n// document order, pairing them up with the replacements list entries\n// that share that find string (also in document order).\nconst uniqueFinds = [...new Set(replacements.map(([find]) => find))];\nconst searchResults = {};\nfor (const find of uniqueFinds) {\n  const res = body.search(find, { matchCase: true, matchWholeWord: false });\n  res.load(\"items\");\n  searchResults[find] = res;\n}\nawait context.sync();\n\nconst nextIndex = {};\nfor (const [find, replace] of replacements) {\n  const idx = nextIndex[find] || 0;\n  const items = searchResults[find].items;\n  if (idx < items.length) {\n    items[idx].insertText(replace, \"Replace\");\n  }\n  nextIndex[find] = idx + 1;\n}\nawait context.sync();\n", "ps1": "# Update the worksheet date and the twenty-five division problems.\n# Listed in document order; \"418\u00f77=\" occurs twice with different\n# replacement values, so each call below uses wdReplaceOne (1) and the\n# Find object's built-in \"search from current position forward\"\n# behaviour to advance past a match it has already replaced.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-06-03 Monday\", \"2024-06-04 Tuesday\"),\n    @(\"274\u00f76=\", \"687\u00f78=\"),\n    @(\"187\u00f76=\", \"305\u00f75=\"),\n    @(\"418\u00f77=\", \"244\u00f72=\"),\n    @(\"700\u00f78=\", \"930\u00f73=\"),\n    @(\"780\u00f77=\", \"846\u00f75=\"),\n    @(\"343\u00f73=\", \"134\u00f73=\"),\n    @(\"336\u00f77=\", \"944\u00f72=\"),\n    @(\"418\u00f77=\", \"708\u00f73=\"),\n    @(\"655\u00f74=\", \"245\u00f79=\"),\n    @(\"290\u00f75=\", \"502\u00f74=\"),\n    @(\"682\u00f76=\", \"923\u00f72=\"),\n    @(\"544\u00f77=\", \"606\u00f79=\"),\n    @(\"590\u00f73=\", \"768\u00f76=\"),\n    @(\"679\u00f74=\", \"160\u00f77=\"),\n    @(\"262\u00f79=\", \"109\u00f74=\"),\n    @(\"821\u00f73=\", \"999\u00f79=\"),\n    @(\"119\u00f75=\", \"755\u00f77=\"),\n    @(\"881\u00f77=\", \"740\u00f78=\"),\n    @(\"652\u00f78=\", \"281\u00f79=\"),\n    @(\"957\u00f73=\", \"183\u00f78=\"),\n    @(\"678\u00f74=\", \"540\u00f75=\"),\n    @(\"395\u00f76=\", \"898\u00f72=\"),\n    @(\"515\u00f78=\", \"768\u00f79=\"),\n    @(\"535\u00f75=\", \"857\u00f76=\"),\n    @(\"625\u00f73=\", \"308\u00f78=\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $find = $d.Content.Find\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n}\n"}
